$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" using the same style as the other header cells (H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data rows 2-21: I column = 1 (constant), J column = copy of H column value
for ($r = 2; $r -le 21; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
